$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$newValues = @(
    "Thu Feb 16 12:14:13 EST 2023",
    "Thu Feb 16 12:14:22 EST 2023",
    "Thu Feb 16 12:14:31 EST 2023",
    "Thu Feb 16 12:14:40 EST 2023",
    "Thu Feb 16 12:14:49 EST 2023",
    "Thu Feb 16 12:14:58 EST 2023",
    "Thu Feb 16 12:15:08 EST 2023",
    "Thu Feb 16 12:15:17 EST 2023",
    "Thu Feb 16 12:15:26 EST 2023",
    "Thu Feb 16 12:15:35 EST 2023",
    "Thu Feb 16 12:15:45 EST 2023",
    "Thu Feb 16 12:15:54 EST 2023",
    "Thu Feb 16 12:16:03 EST 2023",
    "Thu Feb 16 12:16:12 EST 2023",
    "Thu Feb 16 12:16:21 EST 2023"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $newValues[$i]
}
